# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (before the "总计" sheet) holding the
#    fund-holdings detail for the new quarter, built as a structural copy
#    of the previous quarter's sheet so headers/styles/number formats match.
# 2. Prepend a "2022-Q1" summary row to the "总计" (totals) sheet, shifting
#    the existing quarters down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet right before "总计"
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

$template.Copy($totalSheetBefore)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# the template (2021-Q4) had 10 data rows (rows 2-11); the new quarter only
# has 7 (rows 2-8), so drop the trailing 3 rows and shift the rest up
$newSheet.Range("A9:H11").Delete(-4162)

# overwrite the fund-holding details for 2022-Q1
$newSheet.Range("B2:G8").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "217024"
$newSheet.Range("C2").Value = "招商安盈债券"
$newSheet.Range("D2").Value = "35.05"
$newSheet.Range("E2").Value = "20.20"
$newSheet.Range("F2").Value = "1.33"
$newSheet.Range("G2").Value = "0.4662"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "260112"
$newSheet.Range("C3").Value = "景顺长城能源基建混合"
$newSheet.Range("D3").Value = "16.49"
$newSheet.Range("E3").Value = "60.89"
$newSheet.Range("F3").Value = "2.02"
$newSheet.Range("G3").Value = "0.3331"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "014887"
$newSheet.Range("C4").Value = "招商安福1年定期开放债券"
$newSheet.Range("D4").Value = "17.22"
$newSheet.Range("E4").Value = "27.65"
$newSheet.Range("F4").Value = "1.50"
$newSheet.Range("G4").Value = "0.2583"
$newSheet.Range("H4").Value = 6

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "009927"
$newSheet.Range("C5").Value = "工银瑞信聚利18个月定期开放混合A"
$newSheet.Range("D5").Value = "5.54"
$newSheet.Range("E5").Value = "23.27"
$newSheet.Range("F5").Value = "0.96"
$newSheet.Range("G5").Value = "0.0532"
$newSheet.Range("H5").Value = 3

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "000195"
$newSheet.Range("C6").Value = "工银瑞信成长收益混合A"
$newSheet.Range("D6").Value = "4.21"
$newSheet.Range("E6").Value = "31.51"
$newSheet.Range("F6").Value = "0.94"
$newSheet.Range("G6").Value = "0.0396"
$newSheet.Range("H6").Value = 8

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "000196"
$newSheet.Range("C7").Value = "工银瑞信成长收益混合B"
$newSheet.Range("D7").Value = "1.44"
$newSheet.Range("E7").Value = "31.51"
$newSheet.Range("F7").Value = "0.94"
$newSheet.Range("G7").Value = "0.0135"
$newSheet.Range("H7").Value = 8

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "009928"
$newSheet.Range("C8").Value = "工银瑞信聚利18个月定期开放混合C"
$newSheet.Range("D8").Value = "0.83"
$newSheet.Range("E8").Value = "23.27"
$newSheet.Range("F8").Value = "0.96"
$newSheet.Range("G8").Value = "0.0080"
$newSheet.Range("H8").Value = 3

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row into "总计"
# ---------------------------------------------------------------------
# re-fetch by name: sheet references captured before the Copy()/rename
# above track *position*, not identity, so they'd now resolve to the
# newly inserted "2022-Q1" sheet instead of "总计"
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert(-4121)

# re-apply the index-column style (s=2) used by the rest of column A,
# picked up from an untouched sheet's A2 cell
$wb.Worksheets.Item("2021-Q3").Range("A2:D2").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 1.17

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# ---------------------------------------------------------------------
# restore the originally-active sheet/selection
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
$wb.Worksheets.Item("2021-Q1").Range("A1").Select()
